# Scheduled-runner sheet update: refresh market-price derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across the
# per-job Profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
#
# Values below mirror a refreshed market-board pull; rows whose HQ price
# data disappeared have their M/N (profit) cells cleared, and a couple of
# rows gained fresh HQ pricing so their N (LeveProfitHQ) cell is populated
# for the first time.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1580.5625
$ws.Range("I18").Value = 1312.6
$ws.Range("K18").Value = 1312.6
$ws.Range("M18").Value = -1028.6
$ws.Range("H33").Value = 2235.111
$ws.Range("I33").Value = 1722.2307
$ws.Range("J33").Value = 3568.6
$ws.Range("K33").Value = 1722.2307
$ws.Range("L33").Value = 3568.6
$ws.Range("M33").Value = -1493.2307
$ws.Range("N33").Value = -4026.6
$ws.Range("H38").Value = 1147.375
$ws.Range("I38").Value = 1147.375
$ws.Range("K38").Value = 3442.125
$ws.Range("M38").Value = -3070.125
$ws.Range("H64").Value = 4901.75
$ws.Range("I64").Value = 4663.2856
$ws.Range("J64").Value = 4999.9414
$ws.Range("K64").Value = 4663.2856
$ws.Range("L64").Value = 4999.9414
$ws.Range("M64").Value = -4415.2856
$ws.Range("N64").Value = -5495.9414
$ws.Range("H67").Value = 4901.75
$ws.Range("I67").Value = 4663.2856
$ws.Range("J67").Value = 4999.9414
$ws.Range("K67").Value = 4663.2856
$ws.Range("L67").Value = 4999.9414
$ws.Range("M67").Value = -3805.2856
$ws.Range("N67").Value = -6715.9414
$ws.Range("H76").Value = 5948.4546
$ws.Range("I76").Value = 4405.5
$ws.Range("J76").Value = 7800
$ws.Range("K76").Value = 4405.5
$ws.Range("L76").Value = 7800
$ws.Range("M76").Value = -4090.5
$ws.Range("N76").Value = -8430
$ws.Range("H79").Value = 5948.4546
$ws.Range("I79").Value = 4405.5
$ws.Range("J79").Value = 7800
$ws.Range("K79").Value = 4405.5
$ws.Range("L79").Value = 7800
$ws.Range("M79").Value = -3313.5
$ws.Range("N79").Value = -9984
$ws.Range("H82").Value = 20070.875
$ws.Range("I82").Value = 12931.429
$ws.Range("K82").Value = 38794.287
$ws.Range("M82").Value = -38388.287
$ws.Range("H85").Value = 20070.875
$ws.Range("I85").Value = 12931.429
$ws.Range("K85").Value = 38794.287
$ws.Range("M85").Value = -37390.287
$ws.Range("H100").Value = 2732.4783
$ws.Range("I100").Value = 1587
$ws.Range("J100").Value = 3233.625
$ws.Range("K100").Value = 1587
$ws.Range("L100").Value = 3233.625
$ws.Range("M100").Value = -1046
$ws.Range("N100").Value = -4315.625
$ws.Range("H113").Value = 55559776
$ws.Range("I113").Value = 20003996
$ws.Range("J113").Value = 100004500
$ws.Range("K113").Value = 20003996
$ws.Range("L113").Value = 100004500
$ws.Range("M113").Value = -20000742
$ws.Range("N113").Value = -100011008
$ws.Range("H132").Value = 2400.3333
$ws.Range("I132").Value = 1985
$ws.Range("K132").Value = 5955
$ws.Range("M132").Value = -3425
$ws.Range("H135").Value = 2930
$ws.Range("I135").Value = 2757.5
$ws.Range("K135").Value = 24817.5
$ws.Range("M135").Value = -22282.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H122").Value = 3464.625
$ws.Range("I122").Value = 3277
$ws.Range("K122").Value = 9831
$ws.Range("M122").Value = -7381

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H86").Value = 3937.625
$ws.Range("J86").Value = 3632.3333
$ws.Range("L86").Value = 3632.3333
$ws.Range("N86").Value = -5878.3333
$ws.Range("H89").Value = 3937.625
$ws.Range("J89").Value = 3632.3333
$ws.Range("L89").Value = 18161.6665
$ws.Range("N89").Value = -29393.6665
$ws.Range("H99").Value = 3367.5
$ws.Range("I99").Value = 2645.5
$ws.Range("K99").Value = 2645.5
$ws.Range("M99").Value = -1147.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 3700
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H44").Value = 31000
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H46").Value = 3700
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H58").Value = 2204.1765
$ws.Range("I58").Value = 1631.4
$ws.Range("K58").Value = 1631.4
$ws.Range("M58").Value = -1428.4
$ws.Range("H136").Value = 2204.1765
$ws.Range("I136").Value = 1631.4
$ws.Range("K136").Value = 4894.200000000001
$ws.Range("M136").Value = -2344.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 6151.375
$ws.Range("J34").Value = 9600.200000000001
$ws.Range("L34").Value = 28800.6
$ws.Range("N34").Value = -28968.6
$ws.Range("H55").Value = 2009.875
$ws.Range("I55").Value = 867.7143
$ws.Range("J55").Value = 10005
$ws.Range("K55").Value = 2603.1429
$ws.Range("L55").Value = 30015
$ws.Range("M55").Value = -2426.1429
$ws.Range("N55").Value = -30369
$ws.Range("H138").Value = 4407.1113
$ws.Range("I138").Value = 3666
$ws.Range("K138").Value = 10998
$ws.Range("M138").Value = -5858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 186.5
$ws.Range("I2").Value = 58.285713
$ws.Range("J2").Value = 286.22223
$ws.Range("K2").Value = 58.285713
$ws.Range("L2").Value = 286.22223
$ws.Range("M2").Value = 54.714287
$ws.Range("N2").Value = -512.2222300000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1344.5555
$ws.Range("I22").Value = 1266.3334
$ws.Range("K22").Value = 1266.3334
$ws.Range("M22").Value = -971.3334
$ws.Range("H27").Value = 1344.5555
$ws.Range("I27").Value = 1266.3334
$ws.Range("K27").Value = 1266.3334
$ws.Range("M27").Value = -1159.3334
$ws.Range("H38").Value = 38499.5
$ws.Range("J38").Value = 54999
$ws.Range("L38").Value = 54999
$ws.Range("N38").Value = -55819
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H87").Value = 58375
$ws.Range("J87").Value = 67562.5
$ws.Range("L87").Value = 67562.5
$ws.Range("N87").Value = -69808.5
$ws.Range("H90").Value = 58375
$ws.Range("J90").Value = 67562.5
$ws.Range("L90").Value = 202687.5
$ws.Range("N90").Value = -213919.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H96").Value = 5157.6665
$ws.Range("I96").Value = 4989.2
$ws.Range("K96").Value = 4989.2
$ws.Range("M96").Value = -3616.2
$ws.Range("H126").Value = 4100.8
$ws.Range("I126").Value = 4000.6667
$ws.Range("K126").Value = 12002.0001
$ws.Range("M126").Value = -9532.000100000001
$ws.Range("H132").Value = 1820.0476
$ws.Range("I132").Value = 1335
$ws.Range("K132").Value = 4005
$ws.Range("M132").Value = -1475
